$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1831.579
$ws.Range("J112").Value = 2053.2666
$ws.Range("L112").Value = 6159.7998
$ws.Range("N112").Value = -8375.799800000001
$ws.Range("H137").Value = 28573176
$ws.Range("I137").Value = 1284.1786
$ws.Range("J137").Value = 142860740
$ws.Range("K137").Value = 3852.5358
$ws.Range("L137").Value = 428582220
$ws.Range("M137").Value = -1302.5358
$ws.Range("N137").Value = -428587320
$ws.Range("H138").Value = 2929.5115
$ws.Range("I138").Value = 2380.3333
$ws.Range("J138").Value = 3172.59
$ws.Range("K138").Value = 7140.999899999999
$ws.Range("L138").Value = 9517.77
$ws.Range("M138").Value = -2000.999899999999
$ws.Range("N138").Value = -19797.77
$ws.Range("H141").Value = 5566.5454
$ws.Range("I141").Value = 2943.5715
$ws.Range("J141").Value = 7499.263
$ws.Range("K141").Value = 8830.7145
$ws.Range("L141").Value = 22497.789
$ws.Range("M141").Value = -3650.7145
$ws.Range("N141").Value = -32857.789

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 33644.94
$ws.Range("I32").Value = 31562.256
$ws.Range("J32").Value = 36961.816
$ws.Range("K32").Value = 31562.256
$ws.Range("L32").Value = 36961.816
$ws.Range("M32").Value = -31275.256
$ws.Range("N32").Value = -37535.816
$ws.Range("H45").Value = 902.1667
$ws.Range("I45").Value = 902.1667
$ws.Range("K45").Value = 902.1667
$ws.Range("M45").Value = -525.1667
$ws.Range("H122").Value = 1747.4117
$ws.Range("I122").Value = 1731.625
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 5194.875
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -2744.875
$ws.Range("N122").Value = -10900
$ws.Range("H132").Value = 6332.604
$ws.Range("I132").Value = 8046.4546
$ws.Range("J132").Value = 3504.75
$ws.Range("K132").Value = 24139.3638
$ws.Range("L132").Value = 10514.25
$ws.Range("M132").Value = -21609.3638
$ws.Range("N132").Value = -15574.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1926.7142
$ws.Range("I86").Value = 1829.5714
$ws.Range("J86").Value = 2218.1428
$ws.Range("K86").Value = 1829.5714
$ws.Range("L86").Value = 2218.1428
$ws.Range("M86").Value = -706.5714
$ws.Range("N86").Value = -4464.1428
$ws.Range("H89").Value = 1926.7142
$ws.Range("I89").Value = 1829.5714
$ws.Range("J89").Value = 2218.1428
$ws.Range("K89").Value = 9147.857
$ws.Range("L89").Value = 11090.714
$ws.Range("M89").Value = -3531.857
$ws.Range("N89").Value = -22322.714
$ws.Range("H94").Value = 584.875
$ws.Range("I94").Value = 529.1429000000001
$ws.Range("K94").Value = 529.1429000000001
$ws.Range("M94").Value = -78.14290000000005
$ws.Range("H99").Value = 907.3333
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 907.3333
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 907.3333
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -3903.3333
$ws.Range("H128").Value = 1800
$ws.Range("I128").Value = 1800
$ws.Range("K128").Value = 5400
$ws.Range("M128").Value = -2910
$ws.Range("H141").Value = 72000
$ws.Range("J141").Value = 72000
$ws.Range("L141").Value = 72000
$ws.Range("N141").Value = -82360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14929324
$ws.Range("I31").Value = 2275.4644
$ws.Range("J31").Value = 25646180
$ws.Range("K31").Value = 2275.4644
$ws.Range("L31").Value = 25646180
$ws.Range("M31").Value = -1980.4644
$ws.Range("N31").Value = -25646770
$ws.Range("H34").Value = 14929324
$ws.Range("I34").Value = 2275.4644
$ws.Range("J34").Value = 25646180
$ws.Range("K34").Value = 2275.4644
$ws.Range("L34").Value = 25646180
$ws.Range("M34").Value = -2073.4644
$ws.Range("N34").Value = -25646584
$ws.Range("H58").Value = 11496357
$ws.Range("I58").Value = 2189.158
$ws.Range("J58").Value = 33335276
$ws.Range("K58").Value = 2189.158
$ws.Range("L58").Value = 33335276
$ws.Range("M58").Value = -1986.158
$ws.Range("N58").Value = -33335682
$ws.Range("H134").Value = 3310.3096
$ws.Range("I134").Value = 3921.4688
$ws.Range("J134").Value = 1354.6
$ws.Range("K134").Value = 11764.4064
$ws.Range("L134").Value = 4063.8
$ws.Range("M134").Value = -9229.4064
$ws.Range("N134").Value = -9133.799999999999
$ws.Range("H136").Value = 11496357
$ws.Range("I136").Value = 2189.158
$ws.Range("J136").Value = 33335276
$ws.Range("K136").Value = 6567.474
$ws.Range("L136").Value = 100005828
$ws.Range("M136").Value = -4017.474
$ws.Range("N136").Value = -100010928

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1503898.6
$ws.Range("J131").Value = 1588107.2
$ws.Range("L131").Value = 4764321.6
$ws.Range("N131").Value = -4774401.6
$ws.Range("H133").Value = 2542.5
$ws.Range("I133").Value = 2622.6667
$ws.Range("J133").Value = 2408.889
$ws.Range("K133").Value = 7868.000100000001
$ws.Range("L133").Value = 7226.667
$ws.Range("M133").Value = -2808.000100000001
$ws.Range("N133").Value = -17346.667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2800
$ws.Range("I102").Value = 3500
$ws.Range("K102").Value = 3500
$ws.Range("M102").Value = -1878
$ws.Range("H111").Value = 20000
$ws.Range("J111").Value = 20000
$ws.Range("L111").Value = 20000
$ws.Range("N111").Value = -26134
$ws.Range("H122").Value = 18520616
$ws.Range("I122").Value = 30305014
$ws.Range("J122").Value = 2274.5715
$ws.Range("K122").Value = 90915042
$ws.Range("L122").Value = 6823.7145
$ws.Range("M122").Value = -90912592
$ws.Range("N122").Value = -11723.7145
$ws.Range("H126").Value = 7449.5
$ws.Range("I126").Value = 7449.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 22348.5
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -19878.5
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 38464028
$ws.Range("I7").Value = 2238.5
$ws.Range("J7").Value = 55558156
$ws.Range("K7").Value = 2238.5
$ws.Range("L7").Value = 55558156
$ws.Range("M7").Value = -2126.5
$ws.Range("N7").Value = -55558380
$ws.Range("H22").Value = 600
$ws.Range("I22").Value = 701.3333
$ws.Range("J22").Value = 478.4
$ws.Range("K22").Value = 701.3333
$ws.Range("L22").Value = 478.4
$ws.Range("M22").Value = -406.3333
$ws.Range("N22").Value = -1068.4
$ws.Range("H27").Value = 600
$ws.Range("I27").Value = 701.3333
$ws.Range("J27").Value = 478.4
$ws.Range("K27").Value = 701.3333
$ws.Range("L27").Value = 478.4
$ws.Range("M27").Value = -594.3333
$ws.Range("N27").Value = -692.4
$ws.Range("H40").Value = 2724.1667
$ws.Range("I40").Value = 2469
$ws.Range("J40").Value = 4000
$ws.Range("K40").Value = 2469
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -2333
$ws.Range("N40").Value = -4272
$ws.Range("H46").Value = 1364.72
$ws.Range("J46").Value = 1446.2727
$ws.Range("L46").Value = 1446.2727
$ws.Range("N46").Value = -1822.2727
$ws.Range("H61").Value = 14494311
$ws.Range("I61").Value = 1384.0769
$ws.Range("K61").Value = 1384.0769
$ws.Range("M61").Value = -1182.0769
$ws.Range("H82").Value = 2167.75
$ws.Range("I82").Value = 2022
$ws.Range("K82").Value = 2022
$ws.Range("M82").Value = -1661
$ws.Range("H85").Value = 2167.75
$ws.Range("I85").Value = 2022
$ws.Range("K85").Value = 2022
$ws.Range("M85").Value = -774
$ws.Range("H113").Value = 14494311
$ws.Range("I113").Value = 1384.0769
$ws.Range("K113").Value = 1384.0769
$ws.Range("M113").Value = 785.9231
$ws.Range("H126").Value = 38464028
$ws.Range("I126").Value = 2238.5
$ws.Range("J126").Value = 55558156
$ws.Range("K126").Value = 6715.5
$ws.Range("L126").Value = 166674468
$ws.Range("M126").Value = -4245.5
$ws.Range("N126").Value = -166679408

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 12021
$ws.Range("I52").Value = 12021
$ws.Range("K52").Value = 12021
$ws.Range("M52").Value = -11795
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()
